# Generate Report for Handoff
# b.md has been handed off for localization (zh-cn and de-de), so:
#  - Overview sheet: status for b.md in zh-cn / de-de columns -> "Ready for handoff"
#  - zh-cn / de-de sheets: b.md row -> Status "Ready for handoff",
#    new "Latest Handoff File" name/link, new "Latest Handoff Datetime"

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet: row 3 is "b.md"
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# ---------------------------------------------------------------------------
# zh-cn sheet: row 3 is "b.md"
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-10 11:46:38"

# Rebuild the hyperlinks for this sheet so that the "Latest Handoff File"
# hyperlink display text for b.md points at the new handoff file, while
# every other hyperlink on the sheet keeps its previous address/text.
$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1558718923bb52942a1ad5f56c68da32ab496981/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2191bfc761c157cd00b6acb602da70ccbbe939ef/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e0ef38db40778ad2af16896b962e0ba267e3483/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/e2e/b.md", "", "", "b.md")
$zhcn.Hyperlinks.Add($zhcn.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1558718923bb52942a1ad5f56c68da32ab496981/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/2191bfc761c157cd00b6acb602da70ccbbe939ef/e2e/a.md", "", "", "a.md")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/1e0ef38db40778ad2af16896b962e0ba267e3483/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/.localization-config", "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet: row 3 is "b.md"
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-10 11:46:42"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/047e24152cd95612b08b069fc5dbd3d0ef868265/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bf9330d3ed972143a58620ea72d1b165bca49f55/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/44dd44e63b123322faf090d24a00f0d44fbcadf1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/e2e/b.md", "", "", "b.md")
$dede.Hyperlinks.Add($dede.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/047e24152cd95612b08b069fc5dbd3d0ef868265/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/bf9330d3ed972143a58620ea72d1b165bca49f55/e2e/a.md", "", "", "a.md")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/44dd44e63b123322faf090d24a00f0d44fbcadf1/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf", "", "", "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/652e4a485b11d36a86ce1542a5e6e7ac7a242f1d/.localization-config", "", "", ".localization-config")

Write-Host "Handoff report generated."
